$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the new "Time" column (D) ---
$ws.Range("D1").Value = "Time"
$times = @(6, 6, 30, 12, 12, 2, 6, 0)
for ($i = 0; $i -lt $times.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $times[$i]
}

# Match the saved selection state from the authored workbook
$ws.Range("D1:D9").Select()

Write-Host "done"
